$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Slit3"
$ws.Cells.Item(2, 3).Value = "Robo1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"3.590118666666667"
$ws.Cells.Item(2, 8).Value = [double]"10.770356"
$ws.Cells.Item(2, 9).Value = [double]"0.03935865391742773"
$ws.Cells.Item(2, 10).Value = [double]"0.03935865391742773"
$ws.Cells.Item(2, 11).Value = [double]"2"
$ws.Cells.Item(2, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2, 13).Value = [double]"0.4133443333333333"
$ws.Cells.Item(2, 14).Value = [double]"1.240033"
$ws.Cells.Item(2, 15).Value = [double]"0.01404462990513909"
$ws.Cells.Item(2, 16).Value = [double]"0.01404462990513909"
$ws.Cells.Item(2, 17).Value = [double]"1.483955206860889"
$ws.Cells.Item(2, 18).Value = [double]"13.355596861748"
$ws.Cells.Item(2, 19).Value = [double]"0.0005527777278347254"
$ws.Cells.Item(2, 20).Value = [double]"0.0005527777278347254"

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Slit3"
$ws.Cells.Item(3, 3).Value = "Robo1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"3.590118666666667"
$ws.Cells.Item(3, 8).Value = [double]"10.770356"
$ws.Cells.Item(3, 9).Value = [double]"0.03935865391742773"
$ws.Cells.Item(3, 10).Value = [double]"0.03935865391742773"
$ws.Cells.Item(3, 11).Value = [double]"3"
$ws.Cells.Item(3, 12).Value = [double]"1"
$ws.Cells.Item(3, 13).Value = [double]"24.84824866666667"
$ws.Cells.Item(3, 14).Value = [double]"74.544746"
$ws.Cells.Item(3, 15).Value = [double]"0.8442947638833787"
$ws.Cells.Item(3, 16).Value = [double]"0.8442947638833787"
$ws.Cells.Item(3, 17).Value = [double]"89.20816137217511"
$ws.Cells.Item(3, 18).Value = [double]"802.873452349576"
$ws.Cells.Item(3, 19).Value = [double]"0.03323030541598226"
$ws.Cells.Item(3, 20).Value = [double]"0.03323030541598226"

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Slit3"
$ws.Cells.Item(4, 3).Value = "Robo1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"3.590118666666667"
$ws.Cells.Item(4, 8).Value = [double]"10.770356"
$ws.Cells.Item(4, 9).Value = [double]"0.03935865391742773"
$ws.Cells.Item(4, 10).Value = [double]"0.03935865391742773"
$ws.Cells.Item(4, 11).Value = [double]"1"
$ws.Cells.Item(4, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 13).Value = [double]"0.008287000000000001"
$ws.Cells.Item(4, 14).Value = [double]"0.024861"
$ws.Cells.Item(4, 15).Value = [double]"0.0002815760097285016"
$ws.Cells.Item(4, 16).Value = [double]"0.0002815760097285015"
$ws.Cells.Item(4, 17).Value = [double]"0.02975131339066667"
$ws.Cells.Item(4, 18).Value = [double]"0.267761820516"
$ws.Cells.Item(4, 19).Value = [double]"1.108245271835436e-05"
$ws.Cells.Item(4, 20).Value = [double]"1.108245271835436e-05"

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Slit3"
$ws.Cells.Item(5, 3).Value = "Robo1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"3.590118666666667"
$ws.Cells.Item(5, 8).Value = [double]"10.770356"
$ws.Cells.Item(5, 9).Value = [double]"0.03935865391742773"
$ws.Cells.Item(5, 10).Value = [double]"0.03935865391742773"
$ws.Cells.Item(5, 11).Value = [double]"3"
$ws.Cells.Item(5, 12).Value = [double]"1"
$ws.Cells.Item(5, 13).Value = [double]"4.160894333333333"
$ws.Cells.Item(5, 14).Value = [double]"12.482683"
$ws.Cells.Item(5, 15).Value = [double]"0.1413790302017538"
$ws.Cells.Item(5, 16).Value = [double]"0.1413790302017538"
$ws.Cells.Item(5, 17).Value = [double]"14.93810441612755"
$ws.Cells.Item(5, 18).Value = [double]"134.442939745148"
$ws.Cells.Item(5, 19).Value = [double]"0.005564488320892391"
$ws.Cells.Item(5, 20).Value = [double]"0.005564488320892391"

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Slit3"
$ws.Cells.Item(6, 3).Value = "Robo1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = [double]"3"
$ws.Cells.Item(6, 6).Value = [double]"1"
$ws.Cells.Item(6, 7).Value = [double]"76.92488366666667"
$ws.Cells.Item(6, 8).Value = [double]"230.774651"
$ws.Cells.Item(6, 9).Value = [double]"0.8433314202078528"
$ws.Cells.Item(6, 10).Value = [double]"0.8433314202078527"
$ws.Cells.Item(6, 11).Value = [double]"2"
$ws.Cells.Item(6, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6, 13).Value = [double]"0.4133443333333333"
$ws.Cells.Item(6, 14).Value = [double]"1.240033"
$ws.Cells.Item(6, 15).Value = [double]"0.01404462990513909"
$ws.Cells.Item(6, 16).Value = [double]"0.01404462990513909"
$ws.Cells.Item(6, 17).Value = [double]"31.79646475594256"
$ws.Cells.Item(6, 18).Value = [double]"286.168182803483"
$ws.Cells.Item(6, 19).Value = [double]"0.01184427768419463"
$ws.Cells.Item(6, 20).Value = [double]"0.01184427768419463"

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Slit3"
$ws.Cells.Item(7, 3).Value = "Robo1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = [double]"3"
$ws.Cells.Item(7, 6).Value = [double]"1"
$ws.Cells.Item(7, 7).Value = [double]"76.92488366666667"
$ws.Cells.Item(7, 8).Value = [double]"230.774651"
$ws.Cells.Item(7, 9).Value = [double]"0.8433314202078528"
$ws.Cells.Item(7, 10).Value = [double]"0.8433314202078527"
$ws.Cells.Item(7, 11).Value = [double]"3"
$ws.Cells.Item(7, 12).Value = [double]"1"
$ws.Cells.Item(7, 13).Value = [double]"24.84824866666667"
$ws.Cells.Item(7, 14).Value = [double]"74.544746"
$ws.Cells.Item(7, 15).Value = [double]"0.8442947638833787"
$ws.Cells.Item(7, 16).Value = [double]"0.8442947638833787"
$ws.Cells.Item(7, 17).Value = [double]"1911.448638003739"
$ws.Cells.Item(7, 18).Value = [double]"17203.03774203365"
$ws.Cells.Item(7, 19).Value = [double]"0.7120203022998235"
$ws.Cells.Item(7, 20).Value = [double]"0.7120203022998234"

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Slit3"
$ws.Cells.Item(8, 3).Value = "Robo1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = [double]"3"
$ws.Cells.Item(8, 6).Value = [double]"1"
$ws.Cells.Item(8, 7).Value = [double]"76.92488366666667"
$ws.Cells.Item(8, 8).Value = [double]"230.774651"
$ws.Cells.Item(8, 9).Value = [double]"0.8433314202078528"
$ws.Cells.Item(8, 10).Value = [double]"0.8433314202078527"
$ws.Cells.Item(8, 11).Value = [double]"1"
$ws.Cells.Item(8, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 13).Value = [double]"0.008287000000000001"
$ws.Cells.Item(8, 14).Value = [double]"0.024861"
$ws.Cells.Item(8, 15).Value = [double]"0.0002815760097285016"
$ws.Cells.Item(8, 16).Value = [double]"0.0002815760097285015"
$ws.Cells.Item(8, 17).Value = [double]"0.6374765109456668"
$ws.Cells.Item(8, 18).Value = [double]"5.737288598511"
$ws.Cells.Item(8, 19).Value = [double]"0.0002374618961807974"
$ws.Cells.Item(8, 20).Value = [double]"0.0002374618961807974"

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Slit3"
$ws.Cells.Item(9, 3).Value = "Robo1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = [double]"3"
$ws.Cells.Item(9, 6).Value = [double]"1"
$ws.Cells.Item(9, 7).Value = [double]"76.92488366666667"
$ws.Cells.Item(9, 8).Value = [double]"230.774651"
$ws.Cells.Item(9, 9).Value = [double]"0.8433314202078528"
$ws.Cells.Item(9, 10).Value = [double]"0.8433314202078527"
$ws.Cells.Item(9, 11).Value = [double]"3"
$ws.Cells.Item(9, 12).Value = [double]"1"
$ws.Cells.Item(9, 13).Value = [double]"4.160894333333333"
$ws.Cells.Item(9, 14).Value = [double]"12.482683"
$ws.Cells.Item(9, 15).Value = [double]"0.1413790302017538"
$ws.Cells.Item(9, 16).Value = [double]"0.1413790302017538"
$ws.Cells.Item(9, 17).Value = [double]"320.0763125409592"
$ws.Cells.Item(9, 18).Value = [double]"2880.686812868633"
$ws.Cells.Item(9, 19).Value = [double]"0.1192293783276539"
$ws.Cells.Item(9, 20).Value = [double]"0.1192293783276539"

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Slit3"
$ws.Cells.Item(10, 3).Value = "Robo1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = [double]"3"
$ws.Cells.Item(10, 6).Value = [double]"1"
$ws.Cells.Item(10, 7).Value = [double]"0.041643"
$ws.Cells.Item(10, 8).Value = [double]"0.124929"
$ws.Cells.Item(10, 9).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(10, 10).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(10, 11).Value = [double]"2"
$ws.Cells.Item(10, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10, 13).Value = [double]"0.4133443333333333"
$ws.Cells.Item(10, 14).Value = [double]"1.240033"
$ws.Cells.Item(10, 15).Value = [double]"0.01404462990513909"
$ws.Cells.Item(10, 16).Value = [double]"0.01404462990513909"
$ws.Cells.Item(10, 17).Value = [double]"0.017212898073"
$ws.Cells.Item(10, 18).Value = [double]"0.154916082657"
$ws.Cells.Item(10, 19).Value = [double]"6.411855723307977e-06"
$ws.Cells.Item(10, 20).Value = [double]"6.411855723307977e-06"

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Slit3"
$ws.Cells.Item(11, 3).Value = "Robo1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = [double]"3"
$ws.Cells.Item(11, 6).Value = [double]"1"
$ws.Cells.Item(11, 7).Value = [double]"0.041643"
$ws.Cells.Item(11, 8).Value = [double]"0.124929"
$ws.Cells.Item(11, 9).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(11, 10).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(11, 11).Value = [double]"3"
$ws.Cells.Item(11, 12).Value = [double]"1"
$ws.Cells.Item(11, 13).Value = [double]"24.84824866666667"
$ws.Cells.Item(11, 14).Value = [double]"74.544746"
$ws.Cells.Item(11, 15).Value = [double]"0.8442947638833787"
$ws.Cells.Item(11, 16).Value = [double]"0.8442947638833787"
$ws.Cells.Item(11, 17).Value = [double]"1.034755619226"
$ws.Cells.Item(11, 18).Value = [double]"9.312800573034"
$ws.Cells.Item(11, 19).Value = [double]"0.000385449545522288"
$ws.Cells.Item(11, 20).Value = [double]"0.000385449545522288"

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Slit3"
$ws.Cells.Item(12, 3).Value = "Robo1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = [double]"3"
$ws.Cells.Item(12, 6).Value = [double]"1"
$ws.Cells.Item(12, 7).Value = [double]"0.041643"
$ws.Cells.Item(12, 8).Value = [double]"0.124929"
$ws.Cells.Item(12, 9).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(12, 10).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(12, 11).Value = [double]"1"
$ws.Cells.Item(12, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(12, 13).Value = [double]"0.008287000000000001"
$ws.Cells.Item(12, 14).Value = [double]"0.024861"
$ws.Cells.Item(12, 15).Value = [double]"0.0002815760097285016"
$ws.Cells.Item(12, 16).Value = [double]"0.0002815760097285015"
$ws.Cells.Item(12, 17).Value = [double]"0.000345095541"
$ws.Cells.Item(12, 18).Value = [double]"0.003105859869"
$ws.Cells.Item(12, 19).Value = [double]"1.28549115335769e-07"
$ws.Cells.Item(12, 20).Value = [double]"1.28549115335769e-07"

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Slit3"
$ws.Cells.Item(13, 3).Value = "Robo1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = [double]"3"
$ws.Cells.Item(13, 6).Value = [double]"1"
$ws.Cells.Item(13, 7).Value = [double]"0.041643"
$ws.Cells.Item(13, 8).Value = [double]"0.124929"
$ws.Cells.Item(13, 9).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(13, 10).Value = [double]"0.0004565343313861054"
$ws.Cells.Item(13, 11).Value = [double]"3"
$ws.Cells.Item(13, 12).Value = [double]"1"
$ws.Cells.Item(13, 13).Value = [double]"4.160894333333333"
$ws.Cells.Item(13, 14).Value = [double]"12.482683"
$ws.Cells.Item(13, 15).Value = [double]"0.1413790302017538"
$ws.Cells.Item(13, 16).Value = [double]"0.1413790302017538"
$ws.Cells.Item(13, 17).Value = [double]"0.173272122723"
$ws.Cells.Item(13, 18).Value = [double]"1.559449104507"
$ws.Cells.Item(13, 19).Value = [double]"6.454438102517368e-05"
$ws.Cells.Item(13, 20).Value = [double]"6.454438102517368e-05"

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Slit3"
$ws.Cells.Item(14, 3).Value = "Robo1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = [double]"3"
$ws.Cells.Item(14, 6).Value = [double]"1"
$ws.Cells.Item(14, 7).Value = [double]"10.65883866666667"
$ws.Cells.Item(14, 8).Value = [double]"31.976516"
$ws.Cells.Item(14, 9).Value = [double]"0.1168533915433334"
$ws.Cells.Item(14, 10).Value = [double]"0.1168533915433334"
$ws.Cells.Item(14, 11).Value = [double]"2"
$ws.Cells.Item(14, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14, 13).Value = [double]"0.4133443333333333"
$ws.Cells.Item(14, 14).Value = [double]"1.240033"
$ws.Cells.Item(14, 15).Value = [double]"0.01404462990513909"
$ws.Cells.Item(14, 16).Value = [double]"0.01404462990513909"
$ws.Cells.Item(14, 17).Value = [double]"4.405770562780888"
$ws.Cells.Item(14, 18).Value = [double]"39.65193506502799"
$ws.Cells.Item(14, 19).Value = [double]"0.001641162637386428"
$ws.Cells.Item(14, 20).Value = [double]"0.001641162637386428"

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Slit3"
$ws.Cells.Item(15, 3).Value = "Robo1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = [double]"3"
$ws.Cells.Item(15, 6).Value = [double]"1"
$ws.Cells.Item(15, 7).Value = [double]"10.65883866666667"
$ws.Cells.Item(15, 8).Value = [double]"31.976516"
$ws.Cells.Item(15, 9).Value = [double]"0.1168533915433334"
$ws.Cells.Item(15, 10).Value = [double]"0.1168533915433334"
$ws.Cells.Item(15, 11).Value = [double]"3"
$ws.Cells.Item(15, 12).Value = [double]"1"
$ws.Cells.Item(15, 13).Value = [double]"24.84824866666667"
$ws.Cells.Item(15, 14).Value = [double]"74.544746"
$ws.Cells.Item(15, 15).Value = [double]"0.8442947638833787"
$ws.Cells.Item(15, 16).Value = [double]"0.8442947638833787"
$ws.Cells.Item(15, 17).Value = [double]"264.8534736872151"
$ws.Cells.Item(15, 18).Value = [double]"2383.681263184936"
$ws.Cells.Item(15, 19).Value = [double]"0.0986587066220507"
$ws.Cells.Item(15, 20).Value = [double]"0.09865870662205069"

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Slit3"
$ws.Cells.Item(16, 3).Value = "Robo1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = [double]"3"
$ws.Cells.Item(16, 6).Value = [double]"1"
$ws.Cells.Item(16, 7).Value = [double]"10.65883866666667"
$ws.Cells.Item(16, 8).Value = [double]"31.976516"
$ws.Cells.Item(16, 9).Value = [double]"0.1168533915433334"
$ws.Cells.Item(16, 10).Value = [double]"0.1168533915433334"
$ws.Cells.Item(16, 11).Value = [double]"1"
$ws.Cells.Item(16, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16, 13).Value = [double]"0.008287000000000001"
$ws.Cells.Item(16, 14).Value = [double]"0.024861"
$ws.Cells.Item(16, 15).Value = [double]"0.0002815760097285016"
$ws.Cells.Item(16, 16).Value = [double]"0.0002815760097285015"
$ws.Cells.Item(16, 17).Value = [double]"0.08832979603066667"
$ws.Cells.Item(16, 18).Value = [double]"0.7949681642759999"
$ws.Cells.Item(16, 19).Value = [double]"3.290311171401406e-05"
$ws.Cells.Item(16, 20).Value = [double]"3.290311171401405e-05"

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Slit3"
$ws.Cells.Item(17, 3).Value = "Robo1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = [double]"3"
$ws.Cells.Item(17, 6).Value = [double]"1"
$ws.Cells.Item(17, 7).Value = [double]"10.65883866666667"
$ws.Cells.Item(17, 8).Value = [double]"31.976516"
$ws.Cells.Item(17, 9).Value = [double]"0.1168533915433334"
$ws.Cells.Item(17, 10).Value = [double]"0.1168533915433334"
$ws.Cells.Item(17, 11).Value = [double]"3"
$ws.Cells.Item(17, 12).Value = [double]"1"
$ws.Cells.Item(17, 13).Value = [double]"4.160894333333333"
$ws.Cells.Item(17, 14).Value = [double]"12.482683"
$ws.Cells.Item(17, 15).Value = [double]"0.1413790302017538"
$ws.Cells.Item(17, 16).Value = [double]"0.1413790302017538"
$ws.Cells.Item(17, 17).Value = [double]"44.35030140804755"
$ws.Cells.Item(17, 18).Value = [double]"399.1527126724279"
$ws.Cells.Item(17, 19).Value = [double]"0.0165206191721823"
$ws.Cells.Item(17, 20).Value = [double]"0.0165206191721823"
